# Upload analysis from 20200306
#
# The "strain" sheet's row 7 (B7:K7) was mislabeled as the "none_R0_control_mCherry"
# condition; it should read "O2_R0_T_mCherry" instead (matching the notebook's
# actual oxygen/mCherry condition for this row of wells).
$wb = $excel.ActiveWorkbook

$strainSheet = $wb.Worksheets.Item("strain")

# Fix the mislabeled condition across the whole row of wells (B7:K7); the
# outer A/L columns stay as "blank" and are left untouched.
$strainSheet.Range("B7:K7").Value = "O2_R0_T_mCherry"

# The author was last looking at the "strain" sheet with cell J17 selected
# when the workbook was saved (previously "pos_selection" / K11 was the
# active tab).
$strainSheet.Activate()
$strainSheet.Range("J17").Select()
